## "Finish intro and data" — add the mtcars variable dictionary as a new
## worksheet ("Sheet2") after the existing "Sheet1" checklist, make it the
## active sheet, and move the selection on Sheet1 from C11 to B6.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Move the selection on the existing sheet (it is no longer the selected
# tab once Sheet2 is added/activated below, but the cell selection itself
# still needs to move from C11 to B6).
$ws1.Range("B6").Select() | Out-Null

# Add the new sheet right after Sheet1 -> becomes Sheet2 / sheetId 2 and
# is activated (so workbookView.activeTab flips to it and it gets
# tabSelected="1").
$ws2 = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $ws1)

# --- Body rows first (this is the order the strings were typed in, which
# controls the order they were interned into sharedStrings.xml) ---
$ws2.Range("B3").Value = "mpg"
$ws2.Range("C3").Value = "Miles/(US) gallon"

$ws2.Range("B4").Value = "cyl"
$ws2.Range("C4").Value = "Number of cylinders"

$ws2.Range("B5").Value = "disp"
$ws2.Range("C5").Value = "Displacement (cu.in.)"

$ws2.Range("B6").Value = "hp"
$ws2.Range("C6").Value = "Gross horsepower"

$ws2.Range("B7").Value = "drat"
$ws2.Range("C7").Value = "Rear axle ratio"

$ws2.Range("B8").Value = "wt"
$ws2.Range("C8").Value = "Weight (1000 lbs)"

$ws2.Range("B9").Value = "qsec"
$ws2.Range("C9").Value = "1/4 mile time"

$ws2.Range("B10").Value = "vs"
$ws2.Range("C10").Value = "Engine (0 = V-shaped, 1 = straight)"

$ws2.Range("B11").Value = "am"
$ws2.Range("C11").Value = "Transmission (0 = automatic, 1 = manual)"

$ws2.Range("B12").Value = "gear"
$ws2.Range("C12").Value = "Number of forward gears"

$ws2.Range("B13").Value = "carb"
$ws2.Range("C13").Value = "Number of carburetors"

# --- Header row typed in last ---
$ws2.Range("B2").Value = "Variable Name"
$ws2.Range("C2").Value = "Variable Description"
$ws2.Range("D2").Value = "Data Type"

# Column C needs to be wide enough to show the descriptions in full.
$ws2.Columns("C").ColumnWidth = 36.6

# Final view state for the new sheet: zoomed to 145% with D3 selected.
$ws2.Activate() | Out-Null
$excel.ActiveWindow.Zoom = 145
$ws2.Range("D3").Select() | Out-Null

# Keep the page orientation as portrait (matches the author's printed
# page setup for the new sheet).
$ws2.PageSetup.Orientation = 1
